$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.89063024520874
$ws.Range("B1").Value = 3.811607122421265
$ws.Range("C1").Value = 2.258739709854126
$ws.Range("D1").Value = 1.781139254570007
$ws.Range("E1").Value = 1.176329851150513
